$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF" (same style as other headers, e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I and J (rows 2-23)
$values = @(
    @(9, 9),
    @(10, 10),
    @(9, 9),
    @(13, 13),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(6, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
